$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Dr. Eman Tantawi, Dr. Servinaz Sayed Mohammad, Dr. Majorelle Magdy, Dr. Hend Mahmoud"
$ws.Range("G3").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G4").Value = "Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G5").Value = "Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Nesma, Dr. Veronia Rafat, Dr. Hanan Ragab"
$ws.Range("G6").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G7").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Veronia Rafat"
$ws.Range("G8").Value = "Administrator, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G9").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G10").Value = "Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel"
$ws.Range("G11").Value = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G12").Value = "Administrator, Dr. Salma El-Gendy"
$ws.Range("G13").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Shimaa Ashraf, Dr. Omnia Mohammad"
$ws.Range("G17").Value = "Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Marwa Mustafa, Dr. Esraa Mostafa, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Dina Adel"
$ws.Range("G24").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Monica, Dr. Marina Atef, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon, Dr. Aya Emad"
$ws.Range("G25").Value = "Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Eman Samir Gabry, Dr. Aya Emad, Dr. Remon, Dr. Marina Atef, Dr. Abdullah El-Agrody"
$ws.Range("G27").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon"
$ws.Range("G28").Value = "Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G29").Value = "Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon"
$ws.Range("G30").Value = "Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G31").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Amira Sobhy, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G32").Value = "Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Veronia Rafat"
$ws.Range("G33").Value = "Dr. Nourhan Mahmoud, Dr. Mohammad El-Tanany, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Nesma, Dr. Veronia Rafat, Dr. Hanan Ragab"
$ws.Range("G34").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Nahla Nagiub, Dr. Veronia Rafat"
$ws.Range("G35").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Servinaz Sayed Mohammad, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Veronia Rafat"
$ws.Range("G36").Value = "Administrator, Dr. Shimaa Ahmad Mekki, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G37").Value = "Dr. Menna tuâ€™Allah Medhat, Dr. Rana Abo-Zaid, Dr. Amira Sobhy, Dr. Hend Mahmoud, Dr. Asmaa Reda, Dr. Gehan Adel, Dr. Majorelle Magdy, Dr. Manar Montaser"
$ws.Range("G38").Value = "Dr. Sara Wael, Dr. Shimaa Ahmad Mekki, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad, Dr. Gehan Adel"
$ws.Range("G39").Value = "Dr. Veronia Rafat, Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Range("G40").Value = "Administrator, Dr. Salma El-Gendy"
$ws.Range("G41").Value = "Dr. Mariam Nour El-Din, D Wessam Atef, Dr. Safa Hany, Dr. Shimaa Ashraf, Dr. Omnia Mohammad"
$ws.Range("G45").Value = "Dr. Eman M. Abo-Sakaya, Dr. Arwa Al-Sayed, Dr. Marwa Mustafa, Dr. Esraa Mostafa, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen, Dr. Dina Adel"
$ws.Range("G52").Value = "Dr. Ola Abd Al-Fattah, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Youstina Magdy, Dr. Monica, Dr. Marina Atef, Dr. Salma Hassan, Dr. Maryam Ashraf, Dr. Remon, Dr. Aya Emad"
$ws.Range("G53").Value = "Dr. Ola Abd Al-Fattah, Dr. Youstina Magdy, Dr. Eman Samir Gabry, Dr. Aya Emad, Dr. Remon, Dr. Marina Atef, Dr. Abdullah El-Agrody"
$ws.Range("G55").Value = "Dr. Ola Abd Al-Fattah, Dr. Eman Mohammad Al, Dr. Wafaa Ebida, Dr. Yasmin, Dr. Neveen Nashaat, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon"
$ws.Range("G56").Value = "Dr. Nardine, Dr. Wafaa Ebida, Dr. Neveen Nashaat, Dr. Aya Hanafy, Dr. Eman Samir Gabry, Dr. Salma Hassan, Dr. Remon, Dr. Abdullah El-Agrody"
$ws.Range("G57").Value = "Dr. Ola Abd Al-Fattah, Dr. Naema Gomaa, Dr. Neveen Nashaat, Dr. Monica, Dr. Eman Samir Gabry, Dr. Remon"
